# Update the "Export" sheet data:
#  - Remove the rows for accounts 004813166 (VENIA), 004450760 (SILVIO),
#    004488571 (CARLOS), 003921139 (GEISA), the old 004322549 (SIMONE) row,
#    004211368 (ILTON), 004690692 (PHYLIA) and 004261201 (ANA).
#  - Add a new row for 004996634 (HIROKO) ahead of 004376853 (ALBERTO) /
#    005064129 (THIAGO), which stay in place.
#  - Re-add 004322549 (SIMONE) right after THIAGO with an updated balance.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# 1) Delete the six consecutive rows that are removed outright
#    (rows 2-7: VENIA, SILVIO, CARLOS, GEISA, SIMONE[old], ILTON).
$ws.Range("A2:A7").EntireRow.Delete()

# After the delete above, ALBERTO/THIAGO shifted up to rows 2-3, and
# PHYLIA/ANA shifted up to rows 4-5. Delete those two as well.
$ws.Range("A4:A5").EntireRow.Delete()

# Now row 2 = ALBERTO, row 3 = THIAGO, row 4 = GUILHERME (unchanged tail).
# 2) Insert a new row above ALBERTO for HIROKO.
$ws.Range("A2").EntireRow.Insert()
$ws.Cells.Item(2, 1).Value = "'004996634"
$ws.Cells.Item(2, 2).Value = "HIROKO"
$ws.Cells.Item(2, 3).Value = 22533.2

# Row layout is now: 2=HIROKO, 3=ALBERTO, 4=THIAGO, 5=GUILHERME...
# 3) Insert a new row below THIAGO (row 4) for the updated SIMONE balance.
$ws.Range("A5").EntireRow.Insert()
$ws.Cells.Item(5, 1).Value = "'004322549"
$ws.Cells.Item(5, 2).Value = "SIMONE"
$ws.Cells.Item(5, 3).Value = 1153.67
